$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of test-table data (row 5), mirroring the existing row 4 layout.
$ws.Range("B5").Value = "demo.TBItem"
$ws.Range("C5").Value = "demo.Item"
$ws.Range("D5").Value = $true
$ws.Range("E5").Value = "item.xlsx"

# Row 4's B/C use the "s=2" font-only style; the new row instead reuses the
# "好" (Good) style already applied to column headers / E4, so copy that
# formatting across rather than assigning a style by name (which would
# otherwise register a brand-new, duplicate cell style).
$ws.Range("E1").Copy()
$ws.Range("B5:C5").PasteSpecial(-4122)
$ws.Range("E5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reflect the reviewer's zoomed-in view and new selection.
$excel.ActiveWindow.Zoom = 175
$ws.Range("D9").Select() | Out-Null
